$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Date value
$ws.Range("B8").Value = "2025-07-11T12:29:53+00:00"

# Set the Jurisdiction value to FRANCE
$ws.Range("B11").Value = "FRANCE"
